$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

function Copy-Fmt($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ---- row 14 (style pattern from row 2) ----
Copy-Fmt "B2" "B14"
$ws.Range("B14").Value = "Evaluation"
Copy-Fmt "C2" "C14"
$ws.Range("C14").Value = "fixed 4 6-Tupels 200k TDNT2 afterState.agt.zip"
Copy-Fmt "D2" "D14"
Copy-Fmt "E2" "E14"
Copy-Fmt "F2" "F14"

# ---- row 17 (style pattern from row 5) ----
Copy-Fmt "C5" "C17"
$ws.Range("C17").Value = "highest tile"
Copy-Fmt "H5" "H17"
$ws.Range("H17").Value = "sum"
Copy-Fmt "I5" "I17"
$ws.Range("I17").Value = "percent"
Copy-Fmt "J5" "J17"
$ws.Range("J17").Value = "percent cum"

# ---- row 18 (style pattern from row 6) ----
Copy-Fmt "C6" "C18"
$ws.Range("C18").Value = "< 1024"
Copy-Fmt "D6" "D18"
$ws.Range("D18").Value = 2
Copy-Fmt "E6" "E18"
$ws.Range("E18").Value = 2
Copy-Fmt "F6" "F18"
$ws.Range("F18").Value = 1
Copy-Fmt "G6" "G18"
$ws.Range("G18").Value = 0
Copy-Fmt "H6" "H18"
$ws.Range("H18").Formula = "=SUM(D18:G18)"
Copy-Fmt "I6" "I18"
$ws.Range("I18").Formula = "=H18/H`$11"
Copy-Fmt "J6" "J18"
$ws.Range("J18").Formula = "=J19+I18"

# ---- row 19 (style pattern from row 7) ----
Copy-Fmt "C7" "C19"
$ws.Range("C19").Value = 1024
Copy-Fmt "D7" "D19"
$ws.Range("D19").Value = 3
Copy-Fmt "E7" "E19"
$ws.Range("E19").Value = 0
Copy-Fmt "F7" "F19"
$ws.Range("F19").Value = 2
Copy-Fmt "G7" "G19"
$ws.Range("G19").Value = 3
Copy-Fmt "H7" "H19"
$ws.Range("H19").Formula = "=SUM(D19:G19)"
Copy-Fmt "I7" "I19"
$ws.Range("I19").Formula = "=H19/H`$11"
Copy-Fmt "J7" "J19"
$ws.Range("J19").Formula = "=J20+I19"

# ---- row 20 (style pattern from row 8) ----
Copy-Fmt "C8" "C20"
$ws.Range("C20").Value = 2048
Copy-Fmt "D8" "D20"
$ws.Range("D20").Value = 3
Copy-Fmt "E8" "E20"
$ws.Range("E20").Value = 5
Copy-Fmt "F8" "F20"
$ws.Range("F20").Value = 3
Copy-Fmt "G8" "G20"
$ws.Range("G20").Value = 3
Copy-Fmt "H8" "H20"
$ws.Range("H20").Formula = "=SUM(D20:G20)"
Copy-Fmt "I8" "I20"
$ws.Range("I20").Formula = "=H20/H`$11"
Copy-Fmt "J8" "J20"
$ws.Range("J20").Formula = "=J21+I20"

# ---- row 21 (style pattern from row 9) ----
Copy-Fmt "C9" "C21"
$ws.Range("C21").Value = 4096
Copy-Fmt "D9" "D21"
$ws.Range("D21").Value = 17
Copy-Fmt "E9" "E21"
$ws.Range("E21").Value = 27
Copy-Fmt "F9" "F21"
$ws.Range("F21").Value = 23
Copy-Fmt "G9" "G21"
$ws.Range("G21").Value = 13
Copy-Fmt "H9" "H21"
$ws.Range("H21").Formula = "=SUM(D21:G21)"
Copy-Fmt "I9" "I21"
$ws.Range("I21").Formula = "=H21/H`$11"
Copy-Fmt "J9" "J21"
$ws.Range("J21").Formula = "=J22+I21"

# ---- row 22 (style pattern from row 10) ----
Copy-Fmt "C10" "C22"
$ws.Range("C22").Value = 8192
Copy-Fmt "D10" "D22"
$ws.Range("D22").Value = 25
Copy-Fmt "E10" "E22"
$ws.Range("E22").Value = 16
Copy-Fmt "F10" "F22"
$ws.Range("F22").Value = 21
Copy-Fmt "G10" "G22"
$ws.Range("G22").Value = 31
Copy-Fmt "H10" "H22"
$ws.Range("H22").Formula = "=SUM(D22:G22)"
Copy-Fmt "I10" "I22"
$ws.Range("I22").Formula = "=H22/H`$11"
Copy-Fmt "J10" "J22"
$ws.Range("J22").Formula = "=I22"

# ---- row 23 (style pattern from row 11) ----
Copy-Fmt "D11" "D23"
$ws.Range("D23").Formula = "=SUM(D18:D22)"
Copy-Fmt "E11" "E23"
$ws.Range("E23").Formula = "=SUM(E18:E22)"
Copy-Fmt "F11" "F23"
$ws.Range("F23").Formula = "=SUM(F18:F22)"
Copy-Fmt "G11" "G23"
$ws.Range("G23").Formula = "=SUM(G18:G22)"
Copy-Fmt "H11" "H23"
$ws.Range("H23").Formula = "=SUM(H18:H22)"

# ---- row 26 (style pattern from row 2) ----
Copy-Fmt "B2" "B26"
$ws.Range("B26").Value = "Evaluation"
Copy-Fmt "C2" "C26"
$ws.Range("C26").Value = "fixed TEST eTiles 4 6-Tupels 200k TDNT2 afterState.agt.zip"
Copy-Fmt "D2" "D26"
Copy-Fmt "E2" "E26"
Copy-Fmt "F2" "F26"

# ---- row 29 (style pattern from row 5) ----
Copy-Fmt "C5" "C29"
$ws.Range("C29").Value = "highest tile"
Copy-Fmt "H5" "H29"
$ws.Range("H29").Value = "sum"
Copy-Fmt "I5" "I29"
$ws.Range("I29").Value = "percent"
Copy-Fmt "J5" "J29"
$ws.Range("J29").Value = "percent cum"

# ---- row 30 (style pattern from row 6) ----
Copy-Fmt "C6" "C30"
$ws.Range("C30").Value = "< 1024"
Copy-Fmt "D6" "D30"
$ws.Range("D30").Value = 0
Copy-Fmt "E6" "E30"
$ws.Range("E30").Value = 0
Copy-Fmt "F6" "F30"
$ws.Range("F30").Value = 0
Copy-Fmt "G6" "G30"
$ws.Range("G30").Value = 0
Copy-Fmt "H6" "H30"
$ws.Range("H30").Formula = "=SUM(D30:G30)"
Copy-Fmt "I6" "I30"
$ws.Range("I30").Formula = "=H30/H`$11"
Copy-Fmt "J6" "J30"
$ws.Range("J30").Formula = "=J31+I30"

# ---- row 31 (style pattern from row 7) ----
Copy-Fmt "C7" "C31"
$ws.Range("C31").Value = 1024
Copy-Fmt "D7" "D31"
$ws.Range("D31").Value = 1
Copy-Fmt "E7" "E31"
$ws.Range("E31").Value = 1
Copy-Fmt "F7" "F31"
$ws.Range("F31").Value = 1
Copy-Fmt "G7" "G31"
$ws.Range("G31").Value = 0
Copy-Fmt "H7" "H31"
$ws.Range("H31").Formula = "=SUM(D31:G31)"
Copy-Fmt "I7" "I31"
$ws.Range("I31").Formula = "=H31/H`$11"
Copy-Fmt "J7" "J31"
$ws.Range("J31").Formula = "=J32+I31"

# ---- row 32 (style pattern from row 8) ----
Copy-Fmt "C8" "C32"
$ws.Range("C32").Value = 2048
Copy-Fmt "D8" "D32"
$ws.Range("D32").Value = 5
Copy-Fmt "E8" "E32"
$ws.Range("E32").Value = 1
Copy-Fmt "F8" "F32"
$ws.Range("F32").Value = 3
Copy-Fmt "G8" "G32"
$ws.Range("G32").Value = 1
Copy-Fmt "H8" "H32"
$ws.Range("H32").Formula = "=SUM(D32:G32)"
Copy-Fmt "I8" "I32"
$ws.Range("I32").Formula = "=H32/H`$11"
Copy-Fmt "J8" "J32"
$ws.Range("J32").Formula = "=J33+I32"

# ---- row 33 (style pattern from row 9) ----
Copy-Fmt "C9" "C33"
$ws.Range("C33").Value = 4096
Copy-Fmt "D9" "D33"
$ws.Range("D33").Value = 20
Copy-Fmt "E9" "E33"
$ws.Range("E33").Value = 27
Copy-Fmt "F9" "F33"
$ws.Range("F33").Value = 19
Copy-Fmt "G9" "G33"
$ws.Range("G33").Value = 26
Copy-Fmt "H9" "H33"
$ws.Range("H33").Formula = "=SUM(D33:G33)"
Copy-Fmt "I9" "I33"
$ws.Range("I33").Formula = "=H33/H`$11"
Copy-Fmt "J9" "J33"
$ws.Range("J33").Formula = "=J34+I33"

# ---- row 34 (style pattern from row 10) ----
Copy-Fmt "C10" "C34"
$ws.Range("C34").Value = 8192
Copy-Fmt "D10" "D34"
$ws.Range("D34").Value = 24
Copy-Fmt "E10" "E34"
$ws.Range("E34").Value = 21
Copy-Fmt "F10" "F34"
$ws.Range("F34").Value = 27
Copy-Fmt "G10" "G34"
$ws.Range("G34").Value = 23
Copy-Fmt "H10" "H34"
$ws.Range("H34").Formula = "=SUM(D34:G34)"
Copy-Fmt "I10" "I34"
$ws.Range("I34").Formula = "=H34/H`$11"
Copy-Fmt "J10" "J34"
$ws.Range("J34").Formula = "=I34"

# ---- row 35 (style pattern from row 11) ----
Copy-Fmt "D11" "D35"
$ws.Range("D35").Formula = "=SUM(D30:D34)"
Copy-Fmt "E11" "E35"
$ws.Range("E35").Formula = "=SUM(E30:E34)"
Copy-Fmt "F11" "F35"
$ws.Range("F35").Formula = "=SUM(F30:F34)"
Copy-Fmt "G11" "G35"
$ws.Range("G35").Formula = "=SUM(G30:G34)"
Copy-Fmt "H11" "H35"
$ws.Range("H35").Formula = "=SUM(H30:H34)"

